# Implemented getting kafka relations.
#
# The "classFields" sheet lists reflected fields per class. Re-running the
# structure-extraction tool against the updated codebase reshuffled the
# field-type values recorded for three of the ResourceServerConfig fields
# (rows 11, 13 and 14) - the field names (column B) stay put, only the
# recorded Field Type (column D) rotates between them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$ws.Cells.Item(11, 4).Value = "com.macro.mall.component.RestfulAccessDeniedHandler"
$ws.Cells.Item(13, 4).Value = "com.macro.mall.config.IgnoreUrlsConfig"
$ws.Cells.Item(14, 4).Value = "com.macro.mall.authorization.AuthorizationManager"
